$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(28, 8).Value = 852.5714  # H28: 836 -> 852.5714
$ws.Cells.Item(28, 9).Value = 793.8  # I28: 795 -> 793.8
$ws.Cells.Item(28, 10).Value = 999.5  # J28: 1000 -> 999.5
$ws.Cells.Item(28, 11).Value = 793.8  # K28: 795 -> 793.8
$ws.Cells.Item(28, 12).Value = 999.5  # L28: 1000 -> 999.5
$ws.Cells.Item(28, 13).Value = -308.8  # M28: -310 -> -308.8
$ws.Cells.Item(28, 14).Value = -1969.5  # N28: -1970 -> -1969.5
$ws.Cells.Item(62, 8).Value = 6113.7144  # H62: 5600 -> 6113.7144
$ws.Cells.Item(62, 9).Value = 5699.25  # I62: 4333.3335 -> 5699.25
$ws.Cells.Item(62, 10).Value = 6666.3335  # J62: 7500 -> 6666.3335
$ws.Cells.Item(62, 11).Value = 5699.25  # K62: 4333.3335 -> 5699.25
$ws.Cells.Item(62, 12).Value = 6666.3335  # L62: 7500 -> 6666.3335
$ws.Cells.Item(62, 13).Value = -5075.25  # M62: -3709.3335 -> -5075.25
$ws.Cells.Item(62, 14).Value = -7914.3335  # N62: -8748 -> -7914.3335
$ws.Cells.Item(65, 8).Value = 6113.7144  # H65: 5600 -> 6113.7144
$ws.Cells.Item(65, 9).Value = 5699.25  # I65: 4333.3335 -> 5699.25
$ws.Cells.Item(65, 10).Value = 6666.3335  # J65: 7500 -> 6666.3335
$ws.Cells.Item(65, 11).Value = 28496.25  # K65: 21666.6675 -> 28496.25
$ws.Cells.Item(65, 12).Value = 33331.6675  # L65: 37500 -> 33331.6675
$ws.Cells.Item(65, 13).Value = -25376.25  # M65: -18546.6675 -> -25376.25
$ws.Cells.Item(65, 14).Value = -39571.6675  # N65: -43740 -> -39571.6675
$ws.Cells.Item(101, 8).Value = 0  # H101: 300 -> 0
$ws.Cells.Item(101, 9).Value = 0  # I101: 300 -> 0
$ws.Cells.Item(101, 11).Value = 0  # K101: 900 -> 0
$ws.Cells.Item(101, 13).ClearContents()  # M101 was 722
$ws.Cells.Item(113, 8).Value = 0  # H113: 3000 -> 0
$ws.Cells.Item(113, 9).Value = 0  # I113: 3000 -> 0
$ws.Cells.Item(113, 11).Value = 0  # K113: 3000 -> 0
$ws.Cells.Item(113, 13).ClearContents()  # M113 was 254
$ws.Cells.Item(115, 8).Value = 2999.75  # H115: 3599.5 -> 2999.75
$ws.Cells.Item(115, 10).Value = 3000  # J115: 4200 -> 3000
$ws.Cells.Item(115, 12).Value = 9000  # L115: 12600 -> 9000
$ws.Cells.Item(115, 14).Value = -12134  # N115: -15734 -> -12134
$ws.Cells.Item(121, 8).Value = 0  # H121: 300 -> 0
$ws.Cells.Item(121, 10).Value = 0  # J121: 300 -> 0
$ws.Cells.Item(121, 12).ClearContents()  # L121 was 900
$ws.Cells.Item(121, 14).Value = 0  # N121: -4394 -> 0
$ws.Cells.Item(127, 8).Value = 455  # H127: 452.5 -> 455
$ws.Cells.Item(127, 9).Value = 455  # I127: 452.5 -> 455
$ws.Cells.Item(127, 11).Value = 1365  # K127: 1357.5 -> 1365
$ws.Cells.Item(127, 13).Value = 3595  # M127: 3602.5 -> 3595
$ws.Cells.Item(132, 8).Value = 8167.125  # H132: 12041.875 -> 8167.125
$ws.Cells.Item(132, 9).Value = 8167.125  # I132: 12041.875 -> 8167.125
$ws.Cells.Item(132, 11).Value = 24501.375  # K132: 36125.625 -> 24501.375
$ws.Cells.Item(132, 13).Value = -21971.375  # M132: -33595.625 -> -21971.375
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 870.8  # H2: 1042 -> 870.8
$ws.Cells.Item(2, 9).Value = 588.75  # I2: 631.8570999999999 -> 588.75
$ws.Cells.Item(2, 11).Value = 588.75  # K2: 631.8570999999999 -> 588.75
$ws.Cells.Item(2, 13).Value = -475.75  # M2: -518.8570999999999 -> -475.75
$ws.Cells.Item(32, 8).Value = 8045.826  # H32: 8388.817999999999 -> 8045.826
$ws.Cells.Item(32, 9).Value = 3891.9443  # I32: 4091.4707 -> 3891.9443
$ws.Cells.Item(32, 11).Value = 3891.9443  # K32: 4091.4707 -> 3891.9443
$ws.Cells.Item(32, 13).Value = -3604.9443  # M32: -3804.4707 -> -3604.9443
$ws.Cells.Item(61, 8).Value = 16335.833  # H61: 22085.834 -> 16335.833
$ws.Cells.Item(61, 9).Value = 11092.272  # I61: 17365 -> 11092.272
$ws.Cells.Item(61, 11).Value = 11092.272  # K61: 17365 -> 11092.272
$ws.Cells.Item(61, 13).Value = -10880.272  # M61: -17153 -> -10880.272
$ws.Cells.Item(110, 8).Value = 854.875  # H110: 772.125 -> 854.875
$ws.Cells.Item(110, 9).Value = 818.2857  # I110: 772.125 -> 818.2857
$ws.Cells.Item(110, 10).Value = 1111  # J110: 0 -> 1111
$ws.Cells.Item(110, 11).Value = 818.2857  # K110: 772.125 -> 818.2857
$ws.Cells.Item(110, 12).Value = 1111  # L110: 0 -> 1111
$ws.Cells.Item(110, 13).Value = 1226.7143  # M110: 1272.875 -> 1226.7143
$ws.Cells.Item(110, 14).Value = -5201  # N110: None -> -5201
$ws.Cells.Item(116, 8).Value = 870.8  # H116: 1042 -> 870.8
$ws.Cells.Item(116, 9).Value = 588.75  # I116: 631.8570999999999 -> 588.75
$ws.Cells.Item(116, 11).Value = 588.75  # K116: 631.8570999999999 -> 588.75
$ws.Cells.Item(116, 13).Value = 1705.25  # M116: 1662.1429 -> 1705.25
$ws.Cells.Item(130, 8).Value = 17608.666  # H130: 22104.166 -> 17608.666
$ws.Cells.Item(130, 10).Value = 17608.666  # J130: 22104.166 -> 17608.666
$ws.Cells.Item(130, 12).Value = 17608.666  # L130: 22104.166 -> 17608.666
$ws.Cells.Item(130, 14).Value = -27648.666  # N130: -32144.166 -> -27648.666
$ws.Cells.Item(136, 8).Value = 16335.833  # H136: 22085.834 -> 16335.833
$ws.Cells.Item(136, 9).Value = 11092.272  # I136: 17365 -> 11092.272
$ws.Cells.Item(136, 11).Value = 33276.81600000001  # K136: 52095 -> 33276.81600000001
$ws.Cells.Item(136, 13).Value = -30726.81600000001  # M136: -49545 -> -30726.81600000001
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 870.8  # H3: 1042 -> 870.8
$ws.Cells.Item(3, 9).Value = 588.75  # I3: 631.8570999999999 -> 588.75
$ws.Cells.Item(3, 11).Value = 588.75  # K3: 631.8570999999999 -> 588.75
$ws.Cells.Item(3, 13).Value = -474.75  # M3: -517.8570999999999 -> -474.75
$ws.Cells.Item(20, 8).Value = 3000  # H20: 0 -> 3000
$ws.Cells.Item(20, 9).Value = 3000  # I20: 0 -> 3000
$ws.Cells.Item(20, 10).Value = 3000  # J20: 0 -> 3000
$ws.Cells.Item(20, 11).Value = 3000  # K20: 0 -> 3000
$ws.Cells.Item(20, 12).Value = 3000  # L20: 0 -> 3000
$ws.Cells.Item(20, 13).Value = -2753  # M20: None -> -2753
$ws.Cells.Item(20, 14).Value = -3494  # N20: None -> -3494
$ws.Cells.Item(134, 8).Value = 6506  # H134: 1012 -> 6506
$ws.Cells.Item(134, 10).Value = 12000  # J134: 0 -> 12000
$ws.Cells.Item(134, 12).Value = 36000  # L134: 0 -> 36000
$ws.Cells.Item(134, 14).Value = -41070  # N134: None -> -41070
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(86, 8).Value = 1750  # H86: 1372.5 -> 1750
$ws.Cells.Item(86, 9).Value = 0  # I86: 995 -> 0
$ws.Cells.Item(86, 11).Value = 0  # K86: 995 -> 0
$ws.Cells.Item(86, 13).ClearContents()  # M86 was 128
$ws.Cells.Item(89, 8).Value = 1750  # H89: 1372.5 -> 1750
$ws.Cells.Item(89, 9).Value = 0  # I89: 995 -> 0
$ws.Cells.Item(89, 11).Value = 0  # K89: 4975 -> 0
$ws.Cells.Item(89, 13).ClearContents()  # M89 was 641
$ws.Cells.Item(99, 8).Value = 3000  # H99: 900 -> 3000
$ws.Cells.Item(99, 9).Value = 0  # I99: 900 -> 0
$ws.Cells.Item(99, 10).Value = 3000  # J99: 0 -> 3000
$ws.Cells.Item(99, 11).Value = 0  # K99: 900 -> 0
$ws.Cells.Item(99, 12).ClearContents()  # L99 was 0
$ws.Cells.Item(99, 13).Value = 3000  # M99: 598 -> 3000
$ws.Cells.Item(99, 14).Value = -5996  # N99: None -> -5996
$ws.Cells.Item(122, 8).Value = 2011  # H122: 2013.2 -> 2011
$ws.Cells.Item(122, 10).Value = 3332.3333  # J122: 3998.5 -> 3332.3333
$ws.Cells.Item(122, 12).Value = 9996.999899999999  # L122: 11995.5 -> 9996.999899999999
$ws.Cells.Item(122, 14).Value = -14896.9999  # N122: -16895.5 -> -14896.9999
$ws.Cells.Item(126, 8).Value = 3000  # H126: 900 -> 3000
$ws.Cells.Item(126, 9).Value = 0  # I126: 900 -> 0
$ws.Cells.Item(126, 10).Value = 3000  # J126: 0 -> 3000
$ws.Cells.Item(126, 11).Value = 0  # K126: 2700 -> 0
$ws.Cells.Item(126, 12).ClearContents()  # L126 was 0
$ws.Cells.Item(126, 13).Value = 9000  # M126: -230 -> 9000
$ws.Cells.Item(126, 14).Value = -13940  # N126: None -> -13940
$ws.Cells.Item(134, 8).Value = 1974.375  # H134: 2199.2856 -> 1974.375
$ws.Cells.Item(134, 9).Value = 1974.375  # I134: 2199.2856 -> 1974.375
$ws.Cells.Item(134, 11).Value = 5923.125  # K134: 6597.8568 -> 5923.125
$ws.Cells.Item(134, 13).Value = -3388.125  # M134: -4062.8568 -> -3388.125
$ws.Cells.Item(141, 8).Value = 549804.4399999999  # H141: 720295.5600000001 -> 549804.4399999999
$ws.Cells.Item(141, 10).Value = 549804.4399999999  # J141: 720295.5600000001 -> 549804.4399999999
$ws.Cells.Item(141, 12).Value = 549804.4399999999  # L141: 720295.5600000001 -> 549804.4399999999
$ws.Cells.Item(141, 14).Value = -560164.4399999999  # N141: -730655.5600000001 -> -560164.4399999999
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(32, 8).Value = 563.3333  # H32: 550 -> 563.3333
$ws.Cells.Item(32, 9).Value = 600  # I32: 550 -> 600
$ws.Cells.Item(32, 10).Value = 545  # J32: 0 -> 545
$ws.Cells.Item(32, 11).Value = 1800  # K32: 1650 -> 1800
$ws.Cells.Item(32, 12).Value = 1635  # L32: 0 -> 1635
$ws.Cells.Item(32, 13).Value = -1517  # M32: -1367 -> -1517
$ws.Cells.Item(32, 14).Value = -2201  # N32: None -> -2201
$ws.Cells.Item(33, 8).Value = 390  # H33: 270 -> 390
$ws.Cells.Item(33, 10).Value = 0  # J33: 150 -> 0
$ws.Cells.Item(33, 12).Value = 0  # L33: 900 -> 0
$ws.Cells.Item(33, 14).ClearContents()  # N33 was -1466
$ws.Cells.Item(34, 8).Value = 1750.75  # H34: 1105 -> 1750.75
$ws.Cells.Item(34, 9).Value = 500  # I34: 460 -> 500
$ws.Cells.Item(34, 10).Value = 2167.6667  # J34: 1750 -> 2167.6667
$ws.Cells.Item(34, 11).Value = 1500  # K34: 1380 -> 1500
$ws.Cells.Item(34, 12).Value = 6503.000100000001  # L34: 5250 -> 6503.000100000001
$ws.Cells.Item(34, 13).Value = -1416  # M34: -1296 -> -1416
$ws.Cells.Item(34, 14).Value = -6671.000100000001  # N34: -5418 -> -6671.000100000001
$ws.Cells.Item(46, 8).Value = 0  # H46: 25 -> 0
$ws.Cells.Item(46, 9).Value = 0  # I46: 25 -> 0
$ws.Cells.Item(46, 11).Value = 0  # K46: 75 -> 0
$ws.Cells.Item(46, 13).ClearContents()  # M46 was 16
$ws.Cells.Item(63, 8).Value = 3000  # H63: 0 -> 3000
$ws.Cells.Item(63, 10).Value = 3000  # J63: 0 -> 3000
$ws.Cells.Item(63, 12).Value = 9000  # L63: 0 -> 9000
$ws.Cells.Item(63, 14).Value = -10498  # N63: None -> -10498
$ws.Cells.Item(66, 8).Value = 3000  # H66: 0 -> 3000
$ws.Cells.Item(66, 10).Value = 3000  # J66: 0 -> 3000
$ws.Cells.Item(66, 12).Value = 27000  # L66: 0 -> 27000
$ws.Cells.Item(66, 14).Value = -34488  # N66: None -> -34488
$ws.Cells.Item(108, 8).Value = 555.5  # H108: 560.8570999999999 -> 555.5
$ws.Cells.Item(108, 9).Value = 555.5  # I108: 560.8570999999999 -> 555.5
$ws.Cells.Item(108, 11).Value = 1666.5  # K108: 1682.5713 -> 1666.5
$ws.Cells.Item(108, 13).Value = 1213.5  # M108: 1197.4287 -> 1213.5
$ws.Cells.Item(117, 8).Value = 11086.25  # H117: 7622.1665 -> 11086.25
$ws.Cells.Item(117, 9).Value = 0  # I117: 694 -> 0
$ws.Cells.Item(117, 11).Value = 0  # K117: 2082 -> 0
$ws.Cells.Item(117, 13).ClearContents()  # M117 was 1360
$ws.Cells.Item(122, 8).Value = 0  # H122: 4900 -> 0
$ws.Cells.Item(122, 10).Value = 0  # J122: 4900 -> 0
$ws.Cells.Item(122, 12).ClearContents()  # L122 was 44100
$ws.Cells.Item(122, 14).Value = 0  # N122: -49000 -> 0
$ws.Cells.Item(134, 8).Value = 500  # H134: 300 -> 500
$ws.Cells.Item(134, 9).Value = 500  # I134: 300 -> 500
$ws.Cells.Item(134, 11).Value = 1500  # K134: 900 -> 1500
$ws.Cells.Item(134, 13).Value = 3570  # M134: 4170 -> 3570
$ws.Cells.Item(140, 8).Value = 740.625  # H140: 794.3333 -> 740.625
$ws.Cells.Item(140, 9).Value = 740.625  # I140: 794.3333 -> 740.625
$ws.Cells.Item(140, 11).Value = 2221.875  # K140: 2382.9999 -> 2221.875
$ws.Cells.Item(140, 13).Value = 2958.125  # M140: 2797.0001 -> 2958.125
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(47, 8).Value = 95000  # H47: 0 -> 95000
$ws.Cells.Item(47, 10).Value = 95000  # J47: 0 -> 95000
$ws.Cells.Item(47, 12).Value = 95000  # L47: 0 -> 95000
$ws.Cells.Item(47, 14).Value = -96136  # N47: None -> -96136
$ws.Cells.Item(70, 8).Value = 7333  # H70: 6078.143 -> 7333
$ws.Cells.Item(70, 9).Value = 7000  # I70: 5525 -> 7000
$ws.Cells.Item(70, 10).Value = 7499.5  # J70: 6299.4 -> 7499.5
$ws.Cells.Item(70, 11).Value = 7000  # K70: 5525 -> 7000
$ws.Cells.Item(70, 12).Value = 7499.5  # L70: 6299.4 -> 7499.5
$ws.Cells.Item(70, 13).Value = -6730  # M70: -5255 -> -6730
$ws.Cells.Item(70, 14).Value = -8039.5  # N70: -6839.4 -> -8039.5
$ws.Cells.Item(73, 8).Value = 7333  # H73: 6078.143 -> 7333
$ws.Cells.Item(73, 9).Value = 7000  # I73: 5525 -> 7000
$ws.Cells.Item(73, 10).Value = 7499.5  # J73: 6299.4 -> 7499.5
$ws.Cells.Item(73, 11).Value = 7000  # K73: 5525 -> 7000
$ws.Cells.Item(73, 12).Value = 7499.5  # L73: 6299.4 -> 7499.5
$ws.Cells.Item(73, 13).Value = -6064  # M73: -4589 -> -6064
$ws.Cells.Item(73, 14).Value = -9371.5  # N73: -8171.4 -> -9371.5
$ws.Cells.Item(102, 8).Value = 383  # H102: 416 -> 383
$ws.Cells.Item(102, 9).Value = 278.875  # I102: 320.125 -> 278.875
$ws.Cells.Item(102, 11).Value = 278.875  # K102: 320.125 -> 278.875
$ws.Cells.Item(102, 13).Value = 1343.125  # M102: 1301.875 -> 1343.125
$ws.Cells.Item(139, 8).Value = 25000  # H139: 0 -> 25000
$ws.Cells.Item(139, 10).Value = 25000  # J139: 0 -> 25000
$ws.Cells.Item(139, 12).Value = 25000  # L139: 0 -> 25000
$ws.Cells.Item(139, 14).Value = -35280  # N139: None -> -35280
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(32, 8).Value = 1769.25  # H32: 1468 -> 1769.25
$ws.Cells.Item(32, 9).Value = 1769.25  # I32: 1468 -> 1769.25
$ws.Cells.Item(32, 11).Value = 1769.25  # K32: 1468 -> 1769.25
$ws.Cells.Item(32, 13).Value = -1452.25  # M32: -1151 -> -1452.25
$ws.Cells.Item(42, 9).Value = 1000  # I42: 0 -> 1000
$ws.Cells.Item(42, 11).Value = 1000  # K42: 0 -> 1000
$ws.Cells.Item(42, 13).Value = -437  # M42: None -> -437
$ws.Cells.Item(46, 8).Value = 1000  # H46: 0 -> 1000
$ws.Cells.Item(46, 9).Value = 1000  # I46: 0 -> 1000
$ws.Cells.Item(46, 11).Value = 1000  # K46: 0 -> 1000
$ws.Cells.Item(46, 13).Value = -812  # M46: None -> -812
$ws.Cells.Item(49, 9).Value = 1000  # I49: 0 -> 1000
$ws.Cells.Item(49, 11).Value = 1000  # K49: 0 -> 1000
$ws.Cells.Item(49, 13).Value = -853  # M49: None -> -853
$ws.Cells.Item(55, 8).Value = 1056.7858  # H55: 1046.3334 -> 1056.7858
$ws.Cells.Item(55, 10).Value = 1183.6666  # J55: 1143.1428 -> 1183.6666
$ws.Cells.Item(55, 12).Value = 1183.6666  # L55: 1143.1428 -> 1183.6666
$ws.Cells.Item(55, 14).Value = -1529.6666  # N55: -1489.1428 -> -1529.6666
$ws.Cells.Item(141, 8).Value = 100000  # H141: 111111 -> 100000
$ws.Cells.Item(141, 10).Value = 100000  # J141: 111111 -> 100000
$ws.Cells.Item(141, 12).Value = 100000  # L141: 111111 -> 100000
$ws.Cells.Item(141, 14).Value = -110360  # N141: -121471 -> -110360
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(54, 8).Value = 38245.4  # H54: 37757.43 -> 38245.4
$ws.Cells.Item(54, 9).Value = 27999  # I54: 27999.334 -> 27999
$ws.Cells.Item(54, 10).Value = 45076.332  # J54: 45076 -> 45076.332
$ws.Cells.Item(54, 11).Value = 27999  # K54: 27999.334 -> 27999
$ws.Cells.Item(54, 12).Value = 45076.332  # L54: 45076 -> 45076.332
$ws.Cells.Item(54, 13).Value = -27479  # M54: -27479.334 -> -27479
$ws.Cells.Item(54, 14).Value = -46116.332  # N54: -46116 -> -46116.332
$ws.Cells.Item(81, 8).Value = 922.25  # H81: 945 -> 922.25
$ws.Cells.Item(81, 9).Value = 922.25  # I81: 945 -> 922.25
$ws.Cells.Item(81, 11).Value = 1844.5  # K81: 1890 -> 1844.5
$ws.Cells.Item(81, 13).Value = -783.5  # M81: -829 -> -783.5
$ws.Cells.Item(84, 8).Value = 922.25  # H84: 945 -> 922.25
$ws.Cells.Item(84, 9).Value = 922.25  # I84: 945 -> 922.25
$ws.Cells.Item(84, 11).Value = 9222.5  # K84: 9450 -> 9222.5
$ws.Cells.Item(84, 13).Value = -3918.5  # M84: -4146 -> -3918.5
